{"js": "// Replace the 25 \"NN\u00f7N=\" division prompts in the single table with their\n// new values, preserving run/paragraph formatting. Replacements are\n// targeted positionally (table row/column) rather than by searching for\n// the old text, because several old values repeat (e.g. \"27\u00f76=\",\n// \"87\u00f75=\") and some new values collide with other old values\n// (e.g. \"98\u00f77=\", \"36\u00f74=\"), which would make a plain global find/replace\n// ambiguous or order-dependent.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Grid rows 0, 4, 8, 12, 16 hold the five data rows (the others are\n// blank spacer rows); each has 5 columns of \"NN\u00f7N=\" text.\nconst replacements = [\n  // [rowIndex, colIndex, oldText, newText]\n  [0, 0, \"91\u00f76=\", \"28\u00f72=\"],\n  [0, 1, \"36\u00f73=\", \"72\u00f72=\"],\n  [0, 2, \"37\u00f78=\", \"25\u00f72=\"],\n  [0, 3, \"87\u00f75=\", \"70\u00f76=\"],\n  [0, 4, \"90\u00f75=\", \"90\u00f79=\"],\n\n  [4, 0, \"27\u00f76=\", \"41\u00f76=\"],\n  [4, 1, \"11\u00f73=\", \"15\u00f73=\"],\n  [4, 2, \"65\u00f74=\", \"98\u00f77=\"],\n  [4, 3, \"27\u00f76=\", \"94\u00f77=\"],\n  [4, 4, \"47\u00f79=\", \"90\u00f72=\"],\n\n  [8, 0, \"80\u00f76=\", \"68\u00f74=\"],\n  [8, 1, \"21\u00f76=\", \"83\u00f78=\"],\n  [8, 2, \"70\u00f75=\", \"43\u00f74=\"],\n  [8, 3, \"47\u00f74=\", \"36\u00f74=\"],\n  [8, 4, \"17\u00f74=\", \"62\u00f72=\"],\n\n  [12, 0, \"22\u00f76=\", \"77\u00f79=\"],\n  [12, 1, \"11\u00f78=\", \"68\u00f73=\"],\n  [12, 2, \"98\u00f77=\", \"97\u00f75=\"],\n  [12, 3, \"46\u00f76=\", \"49\u00f77=\"],\n  [12, 4, \"67\u00f77=\", \"87\u00f79=\"],\n\n  [16, 0, \"43\u00f79=\", \"15\u00f72=\"],\n  [16, 1, \"65\u00f75=\", \"56\u00f75=\"],\n  [16, 2, \"36\u00f74=\", \"32\u00f73=\"],\n  [16, 3, \"68\u00f78=\", \"27\u00f77=\"],\n  [16, 4, \"87\u00f75=\", \"68\u00f75=\"],\n];\n\n// Grab the cells up front and verify their current text matches what we\n// expect before mutating anything.\nconst cells = replacements.map(([row, col]) => table.getCell(row, col));\nfor (const cell of cells) {\n  cell.load(\"value\");\n}\nawait context.sync();\n\nreplacements.forEach(([row, col, oldText], i) => {\n  const actual = cells[i].value;\n  if (actual !== oldText) {\n    throw new Error(\n      `Cell (${row},${col}) expected \"${oldText}\" but found \"${actual}\"`\n    );\n  }\n});\n\nreplacements.forEach(([, , , newText], i) => {\n  const range = cells[i].body.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n});\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"NN\u00f7N=\" division prompts in the single table with their\n# new values, preserving run/paragraph formatting. Replacements are\n# targeted positionally (table row/column) rather than by searching for\n# the old text, because several old values repeat (e.g. \"27\u00f76=\",\n# \"87\u00f75=\") and some new values collide with other old values\n# (e.g. \"98\u00f77=\", \"36\u00f74=\"), which would make a plain global find/replace\n# ambiguous or order-dependent.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# Word table rows are 1-indexed; rows 1, 5, 9, 13, 17 hold the five data\n# rows (the rows in between are blank spacer rows); each has 5 columns\n# (1..5) of \"NN\u00f7N=\" text.\n$replacements = @(\n    @{ Row = 1; Col = 1; Old = \"91\u00f76=\"; Text = \"28\u00f72=\" },\n    @{ Row = 1; Col = 2; Old = \"36\u00f73=\"; Text = \"72\u00f72=\" },\n    @{ Row = 1; Col = 3; Old = \"37\u00f78=\"; Text = \"25\u00f72=\" },\n    @{ Row = 1; Col = 4; Old = \"87\u00f75=\"; Text = \"70\u00f76=\" },\n    @{ Row = 1; Col = 5; Old = \"90\u00f75=\"; Text = \"90\u00f79=\" },\n\n    @{ Row = 5; Col = 1; Old = \"27\u00f76=\"; Text = \"41\u00f76=\" },\n    @{ Row = 5; Col = 2; Old = \"11\u00f73=\"; Text = \"15\u00f73=\" },\n    @{ Row = 5; Col = 3; Old = \"65\u00f74=\"; Text = \"98\u00f77=\" },\n    @{ Row = 5; Col = 4; Old = \"27\u00f76=\"; Text = \"94\u00f77=\" },\n    @{ Row = 5; Col = 5; Old = \"47\u00f79=\"; Text = \"90\u00f72=\" },\n\n    @{ Row = 9; Col = 1; Old = \"80\u00f76=\"; Text = \"68\u00f74=\" },\n    @{ Row = 9; Col = 2; Old = \"21\u00f76=\"; Text = \"83\u00f78=\" },\n    @{ Row = 9; Col = 3; Old = \"70\u00f75=\"; Text = \"43\u00f74=\" },\n    @{ Row = 9; Col = 4; Old = \"47\u00f74=\"; Text = \"36\u00f74=\" },\n    @{ Row = 9; Col = 5; Old = \"17\u00f74=\"; Text = \"62\u00f72=\" },\n\n    @{ Row = 13; Col = 1; Old = \"22\u00f76=\"; Text = \"77\u00f79=\" },\n    @{ Row = 13; Col = 2; Old = \"11\u00f78=\"; Text = \"68\u00f73=\" },\n    @{ Row = 13; Col = 3; Old = \"98\u00f77=\"; Text = \"97\u00f75=\" },\n    @{ Row = 13; Col = 4; Old = \"46\u00f76=\"; Text = \"49\u00f77=\" },\n    @{ Row = 13; Col = 5; Old = \"67\u00f77=\"; Text = \"87\u00f79=\" },\n\n    @{ Row = 17; Col = 1; Old = \"43\u00f79=\"; Text = \"15\u00f72=\" },\n    @{ Row = 17; Col = 2; Old = \"65\u00f75=\"; Text = \"56\u00f75=\" },\n    @{ Row = 17; Col = 3; Old = \"36\u00f74=\"; Text = \"32\u00f73=\" },\n    @{ Row = 17; Col = 4; Old = \"68\u00f78=\"; Text = \"27\u00f77=\" },\n    @{ Row = 17; Col = 5; Old = \"87\u00f75=\"; Text = \"68\u00f75=\" }\n)\n\n# Verify the current text of every target cell before mutating anything,\n# so a structural mismatch fails loudly instead of silently mis-editing.\nforeach ($rep in $replacements) {\n    $cell = $tbl.Cell($rep.Row, $rep.Col)\n    $actual = $cell.Range.Text.TrimEnd([char]7, [char]13, [char]10)\n    if ($actual -ne $rep.Old) {\n        throw \"Cell ($($rep.Row),$($rep.Col)) expected '$($rep.Old)' but found '$actual'\"\n    }\n}\n\nforeach ($rep in $replacements) {\n    $cell = $tbl.Cell($rep.Row, $rep.Col)\n    $cell.Range.Text = $rep.Text\n}\n"}
